$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.413.42'
$ws.Range('E2').Value = '  +0.11%  '
$ws.Range('D3').Value = '3.676.45'
$ws.Range('E3').Value = '  -0.28%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '643.29'
$ws.Range('D5').ClearFormats()
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '159.93'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.56%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.499'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.79%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.10'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.49%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.448'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.06%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000233'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.55%  '
$ws.Range('D13').Value = '4.294.95'
$ws.Range('E13').Value = '  -0.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.73'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.99%  '
$ws.Range('D15').Value = '3.673.20'
$ws.Range('E15').Value = '  -0.55%  '
$ws.Range('D16').Value = '69.381.03'
$ws.Range('E16').Value = '  +0.10%  '
$ws.Range('E17').Value = '  +0.21%  '
$ws.Range('E18').Value = '  +0.23%  '
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '465.96'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.41%  '
$ws.Range('E21').Value = '  +0.44%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '79.48'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.63%  '
$ws.Range('D24').Value = '3.821.94'
$ws.Range('E24').Value = '  -0.30%  '
$ws.Range('E25').Value = '  +0.20%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000127'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +3.43%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.93'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.22%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.09'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.62%  '
$ws.Range('E29').Value = '  -2.79%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.73'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.66%  '
$ws.Range('E31').Value = '  +0.82%  '
$ws.Range('E32').Value = '  -0.59%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '26.90'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.09%  '
$ws.Range('E34').Value = '  +4.10%  '
$ws.Range('E35').Value = '  -1.72%  '
$ws.Range('D36').Value = '3.668.85'
$ws.Range('E36').Value = '  -0.22%  '
$ws.Range('E37').Value = '  +1.36%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.89'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -5.99%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '179.07'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +4.92%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('D41').ClearFormats()
$ws.Range('E42').Value = '  -0.36%  '
$ws.Range('E43').Value = '  -1.72%  '
$ws.Range('E44').Value = '  -1.68%  '
$ws.Range('E45').Value = '  -1.69%  '
$ws.Range('E46').Value = '  +2.57%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '27.49'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.57%  '
$ws.Range('E48').Value = '  -1.71%  '
$ws.Range('E49').Value = '  -3.24%  '
$ws.Range('E51').Value = '  -3.78%  '
